# Add 2022-Q4 data:
#  - insert a new worksheet "2022-Q4" right after "总计" with the quarter's
#    fund holdings, shifting the other quarter sheets along (they keep
#    their own names/content, just move one slot to the right)
#  - insert a new summary row in "总计" for the 2022-Q4 aggregate stats

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet by duplicating "2022-Q1", which
#    already has the same 2-data-row / A1:H3 shape and carries the right
#    cell styles (header row + index column use the bold/bordered style).
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q1")
$templateSheet.Copy($null, $totalSheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Row 2: 010965 / 中银鑫新消费成长混合A
$q4Sheet.Range("B2:G2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "010965"
$q4Sheet.Range("C2").Value = "中银鑫新消费成长混合A"
$q4Sheet.Range("D2").Value = "3.24"
$q4Sheet.Range("E2").Value = "88.64"
$q4Sheet.Range("F2").Value = "2.49"
$q4Sheet.Range("G2").Value = "0.0807"
$q4Sheet.Range("H2").Value = 9

# Row 3: 010962 / 中银鑫新消费成长混合C
$q4Sheet.Range("B3:G3").NumberFormat = "@"
$q4Sheet.Range("B3").Value = "010962"
$q4Sheet.Range("C3").Value = "中银鑫新消费成长混合C"
$q4Sheet.Range("D3").Value = "0.69"
$q4Sheet.Range("E3").Value = "88.64"
$q4Sheet.Range("F3").Value = "2.49"
$q4Sheet.Range("G3").Value = "0.0172"
$q4Sheet.Range("H3").Value = 9

# The "@" number format forced the literal text (instead of Excel parsing it
# back into a number) - clear the residual format now that the values are in
# place so the cells end up styleless again, matching the other quarter sheets.
$q4Sheet.Range("B2:G3").ClearFormats()

# ---------------------------------------------------------------------
# 2. Insert a new row 2 in "总计" for the 2022-Q4 summary and fix up the
#    formatting (the row-insert inherits the header's bold/border style
#    for B:D, which the data rows don't actually use; only column A does).
# ---------------------------------------------------------------------
$totalSheet.Rows("2:2").Insert()

$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# The physical row-insert shifted the old rows 2-6 down to 3-7 but left their
# index-column (A) values untouched (0,1,2,3,4); the index is really a
# 0-based row-position counter, so rows 3-7 need to become 1,2,3,4,5.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.1
